$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 2.1
$ws.Range("X2").Value = 8.5
$ws.Range("AG2").Value = 351
$ws.Range("AL2").Value = 34
$ws.Range("AO2").Value = 11
$ws.Range("BA2").Value = 101
$ws.Range("BB2").Value = 251
$ws.Range("G4").Value = 1.62
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 6.25
$ws.Range("K4").Value = 2.1
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("X4").Value = 6.5
$ws.Range("AF4").Value = 81
$ws.Range("AI4").Value = 29
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 8.5
$ws.Range("AV4").Value = 81
$ws.Range("G5").Value = 1.36
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 9.75
$ws.Range("J5").Value = 1.85
$ws.Range("K5").Value = 2.18
$ws.Range("L5").Value = 8.5
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 9.800000000000001
$ws.Range("O5").Value = 1.31
$ws.Range("P5").Value = 2.9
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.75
$ws.Range("T5").Value = 2.55
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.52
$ws.Range("W5").Value = 5.3
$ws.Range("X5").Value = 5.5
$ws.Range("Y5").Value = 8.5
$ws.Range("Z5").Value = 8.25
$ws.Range("AA5").Value = 12.5
$ws.Range("AB5").Value = 37
$ws.Range("AC5").Value = 8.5
$ws.Range("AD5").Value = 8.25
$ws.Range("AE5").Value = 25
$ws.Range("AH5").Value = 20
$ws.Range("AI5").Value = 75
$ws.Range("AJ5").Value = 32
$ws.Range("AK5").Value = 350
$ws.Range("AL5").Value = 175
$ws.Range("AM5").Value = 120
$ws.Range("AN5").Value = 2.95
$ws.Range("AO5").Value = 6.1
$ws.Range("AP5").Value = 18.5
$ws.Range("AQ5").Value = 17.5
$ws.Range("AU5").Value = 9
$ws.Range("AV5").Value = 100
$ws.Range("AW5").Value = 10
$ws.Range("AX5").Value = 65
$ws.Range("AY5").Value = 60
$ws.Range("BA5").Value = 500
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 14.8
$ws.Range("S6").Value = 1.38
$ws.Range("T6").Value = 3.06
